$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (D) and volume-change (E) values.
# D-column cells that hold plain decimal numbers are forced to Text
# number format first so Excel does not silently coerce them to
# numeric values (which would, e.g., turn "9.20" into 9.2).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.944.11"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.553.12"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.83"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.16"
$ws.Range("E8").Value = "  +3.92%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.775.09"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.554.18"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.948.04"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.66"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.57"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0697"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.29"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.20"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.35"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.424.54"
$ws.Range("E33").Value = "  +4.65%  "
$ws.Range("E34").Value = "  +4.33%  "
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.976"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.520"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("E41").Value = "  +3.66%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +4.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.991"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.23"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.688.90"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.62"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0952"
$ws.Range("E51").Value = "  +0.54%  "
